$wb = $excel.ActiveWorkbook

# --- Update "Hoja1" A1 conversion text ---
$wsHoja1 = $wb.Worksheets.Item("Hoja1")
$wsHoja1.Range("A1").Value = "Conversión del día 💰`n✅ Dólar paralelo: 68`n`nBinance`n✅ 1000 Bs = 1.76 = 6397.3 pesos`n✅ 6397.3 pesos = 1.75 = 922.85 Bs`n`nPromedio competencia`n✅ Tasa pesos: 20`n✅ Tasa Bs: 20`n✅ % Ganancia: 20%"

# --- Update "tasas" sheet N10/O10/N12/O12 ---
$wsTasas = $wb.Worksheets.Item("tasas")
$wsTasas.Range("N10").Value = 568.99
$wsTasas.Range("O10").Value = 3640
$wsTasas.Range("N12").Value = 3651.99
$wsTasas.Range("O12").Value = 526.822
